$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = '2025-12-25 00:50:39'
$ws.Cells.Item(7, 2).Value = 'Admin'
$ws.Cells.Item(7, 3).Value = 'Login'
$ws.Cells.Item(7, 4).Value = 'login_success'
$ws.Cells.Item(7, 5).Value = 'Role: admin'

$ws.Cells.Item(8, 1).Value = '2025-12-25 00:50:39'
$ws.Cells.Item(8, 2).Value = 'Admin'
$ws.Cells.Item(8, 3).Value = 'dashboard'
$ws.Cells.Item(8, 4).Value = 'access_granted'
$ws.Cells.Item(8, 5).Value = 'Opened dashboard page'

$ws.Cells.Item(9, 1).Value = '2025-12-25 00:50:43'
$ws.Cells.Item(9, 2).Value = 'Admin'
$ws.Cells.Item(9, 3).Value = 'invoice'
$ws.Cells.Item(9, 4).Value = 'access_granted'
$ws.Cells.Item(9, 5).Value = 'Opened invoice page'

$ws.Cells.Item(10, 1).Value = '2025-12-25 00:50:52'
$ws.Cells.Item(10, 2).Value = 'Admin'
$ws.Cells.Item(10, 3).Value = 'invoice'
$ws.Cells.Item(10, 4).Value = 'access_granted'
$ws.Cells.Item(10, 5).Value = 'Opened invoice page'

$ws.Cells.Item(11, 1).Value = '2025-12-25 00:50:52'
$ws.Cells.Item(11, 2).Value = 'Admin'
$ws.Cells.Item(11, 3).Value = 'invoice'
$ws.Cells.Item(11, 4).Value = 'access_granted'
$ws.Cells.Item(11, 5).Value = 'Opened invoice page'

$ws.Cells.Item(12, 1).Value = '2025-12-25 00:50:55'
$ws.Cells.Item(12, 2).Value = 'Admin'
$ws.Cells.Item(12, 3).Value = 'invoice'
$ws.Cells.Item(12, 4).Value = 'access_granted'
$ws.Cells.Item(12, 5).Value = 'Opened invoice page'

$ws.Cells.Item(13, 1).Value = '2025-12-25 00:51:00'
$ws.Cells.Item(13, 2).Value = 'Admin'
$ws.Cells.Item(13, 3).Value = 'invoice'
$ws.Cells.Item(13, 4).Value = 'access_granted'
$ws.Cells.Item(13, 5).Value = 'Opened invoice page'

$ws.Cells.Item(14, 1).Value = '2025-12-25 00:51:04'
$ws.Cells.Item(14, 2).Value = 'Admin'
$ws.Cells.Item(14, 3).Value = 'invoice'
$ws.Cells.Item(14, 4).Value = 'access_granted'
$ws.Cells.Item(14, 5).Value = 'Opened invoice page'

$ws.Cells.Item(15, 1).Value = '2025-12-25 00:51:07'
$ws.Cells.Item(15, 2).Value = 'Admin'
$ws.Cells.Item(15, 3).Value = 'invoice'
$ws.Cells.Item(15, 4).Value = 'access_granted'
$ws.Cells.Item(15, 5).Value = 'Opened invoice page'

$ws.Cells.Item(16, 1).Value = '2025-12-25 00:51:16'
$ws.Cells.Item(16, 2).Value = 'Admin'
$ws.Cells.Item(16, 3).Value = 'invoice'
$ws.Cells.Item(16, 4).Value = 'access_granted'
$ws.Cells.Item(16, 5).Value = 'Opened invoice page'

$ws.Cells.Item(17, 1).Value = '2025-12-25 00:51:18'
$ws.Cells.Item(17, 2).Value = 'Admin'
$ws.Cells.Item(17, 3).Value = 'invoice'
$ws.Cells.Item(17, 4).Value = 'access_granted'
$ws.Cells.Item(17, 5).Value = 'Opened invoice page'

$ws.Cells.Item(18, 1).Value = '2025-12-25 00:52:14'
$ws.Cells.Item(18, 2).Value = 'Admin'
$ws.Cells.Item(18, 3).Value = 'invoice'
$ws.Cells.Item(18, 4).Value = 'access_granted'
$ws.Cells.Item(18, 5).Value = 'Opened invoice page'

$ws.Cells.Item(19, 1).Value = '2025-12-25 00:52:14'
$ws.Cells.Item(19, 2).Value = 'Admin'
$ws.Cells.Item(19, 3).Value = 'invoice'
$ws.Cells.Item(19, 4).Value = 'access_granted'
$ws.Cells.Item(19, 5).Value = 'Opened invoice page'

$ws.Cells.Item(20, 1).Value = '2025-12-25 00:52:15'
$ws.Cells.Item(20, 2).Value = 'Admin'
$ws.Cells.Item(20, 3).Value = 'invoice'
$ws.Cells.Item(20, 4).Value = 'access_granted'
$ws.Cells.Item(20, 5).Value = 'Opened invoice page'

$ws.Cells.Item(21, 1).Value = '2025-12-25 00:52:16'
$ws.Cells.Item(21, 2).Value = 'Admin'
$ws.Cells.Item(21, 3).Value = 'invoice'
$ws.Cells.Item(21, 4).Value = 'access_granted'
$ws.Cells.Item(21, 5).Value = 'Opened invoice page'

$ws.Cells.Item(22, 1).Value = '2025-12-25 00:52:17'
$ws.Cells.Item(22, 2).Value = 'Admin'
$ws.Cells.Item(22, 3).Value = 'invoice'
$ws.Cells.Item(22, 4).Value = 'access_granted'
$ws.Cells.Item(22, 5).Value = 'Opened invoice page'

$ws.Cells.Item(23, 1).Value = '2025-12-25 00:52:17'
$ws.Cells.Item(23, 2).Value = 'Admin'
$ws.Cells.Item(23, 3).Value = 'invoice'
$ws.Cells.Item(23, 4).Value = 'access_granted'
$ws.Cells.Item(23, 5).Value = 'Opened invoice page'

$ws.Cells.Item(24, 1).Value = '2025-12-25 00:56:35'
$ws.Cells.Item(24, 2).Value = 'Admin'
$ws.Cells.Item(24, 3).Value = 'Login'
$ws.Cells.Item(24, 4).Value = 'login_success'
$ws.Cells.Item(24, 5).Value = 'Role: admin'

$ws.Cells.Item(25, 1).Value = '2025-12-25 00:56:35'
$ws.Cells.Item(25, 2).Value = 'Admin'
$ws.Cells.Item(25, 3).Value = 'dashboard'
$ws.Cells.Item(25, 4).Value = 'access_granted'
$ws.Cells.Item(25, 5).Value = 'Opened dashboard page'

$ws.Cells.Item(26, 1).Value = '2025-12-25 00:56:36'
$ws.Cells.Item(26, 2).Value = 'Admin'
$ws.Cells.Item(26, 3).Value = 'quotation'
$ws.Cells.Item(26, 4).Value = 'access_granted'
$ws.Cells.Item(26, 5).Value = 'Opened quotation page'
